# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet
#    that carries that status cell (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2) The two "Latest Handoff/Handback" status columns on Overview (E,F) and
#    the "Status" column (C) on the zh-cn/de-de sheets get narrower, from
#    roughly 17.22 chars to roughly 13.41 chars (matches the shorter text).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 12.5

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "In Translation"
$ws.Columns.Item(3).ColumnWidth = 12.5

$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "In Translation"
$ws.Columns.Item(3).ColumnWidth = 12.5
